$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 - Team Meeting: add Wednesday (D8) hours, update Weekly Total (I8)
$ws.Range("D8").Value = 0.5
$ws.Range("I8").Value = 0.5

# Row 9 - Mentor Meeting: add Wednesday (D9) hours, update Weekly Total (I9)
$ws.Range("D9").Value = 0.5
$ws.Range("I9").Value = 0.5

# Row 13 - Sponsor Work: add Thursday (E13) hours, update Weekly Total (I13)
$ws.Range("E13").Value = 1
$ws.Range("I13").Value = 5

# Row 14 - Daily Total: update Wednesday (D14), Thursday (E14), Weekly Total (I14)
$ws.Range("D14").Value = 3
$ws.Range("E14").Value = 1
$ws.Range("I14").Value = 6

# Update selection to K11
$ws.Range("K11").Select()
